# Applies the rubric update for Milestone 1 (normal mapping on a cube completed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark additional rubric rows as completed for Milestone I ("I" in column E, "X" in column F).
$ws.Range("E14").Value = "I"
$ws.Range("F14").Value = "X"
$ws.Range("F22").Value = "X"
$ws.Range("F23").Value = "X"
$ws.Range("F24").Value = "X"
$ws.Range("F40").Value = "X"

# Recalculate formulas so dependent totals (G/H/I/J/K/L columns) update.
$excel.CalculateFullRebuild()

# Restore the view to the top of the sheet with F16 selected, matching the saved state.
$ws.Activate()
$ws.Range("F16").Select()
